# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (E,K: 1 -> 3) and recompute
# the dependent expression / specificity metrics for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3
        G = 2.520808
        H = 7.562424
        K = 3
        M = 0.9623583333333334
        N = 2.887075
        O = 0.05805926999654511
        P = 0.05805926999654509
        Q = 2.425920585533334
        R = 21.8332852698
        S = 0.05805926999654511
        T = 0.05805926999654509
    }
    3 = @{
        E = 3
        G = 2.520808
        H = 7.562424
        K = 3
        M = 12.12890133333333
        N = 36.386704
        O = 0.731739034081334
        P = 0.7317390340813339
        Q = 30.57463151227734
        R = 275.171683610496
        S = 0.731739034081334
        T = 0.7317390340813339
    }
    4 = @{
        E = 3
        G = 2.520808
        H = 7.562424
        K = 3
        M = 3.484187
        N = 10.452561
        O = 0.210201695922121
        P = 0.2102016959221209
        Q = 8.782966463096001
        R = 79.04669816786399
        S = 0.210201695922121
        T = 0.2102016959221209
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
